$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 4-5, pushing the existing rows 4 and 5
# (29-Jun-2021 / 10-Aug-2021 entries) down to rows 6 and 7.
$ws.Rows("4:5").Insert()

# New row 4: Mapocho Venta Directa de Santiago - Alcachofa, week of 27-Aug-2021
$ws.Range("A4").Value = 12
$ws.Range("B4").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C4").Value = "Metropolitana"
$ws.Range("D4").Value = 44435
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = 100112013
$ws.Range("G4").Value = "Alcachofa"
$ws.Range("H4").Value = "Española"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 14000
$ws.Range("N4").Value = "`$/caja 30 unidades"
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 467
$ws.Range("Q4").Value = 30
$ws.Range("R4").Value = "Hortaliza"

# New row 5: Mapocho Venta Directa de Santiago - Alcachofa, week of 27-Aug-2021 (Elquí origin)
$ws.Range("A5").Value = 12
$ws.Range("B5").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C5").Value = "Metropolitana"
$ws.Range("D5").Value = 44435
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = 100112013
$ws.Range("G5").Value = "Alcachofa"
$ws.Range("H5").Value = "Española"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 25
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 14000
$ws.Range("M5").Value = 14000
$ws.Range("N5").Value = "`$/caja 30 unidades"
$ws.Range("O5").Value = "Provincia del Elquí"
$ws.Range("P5").Value = 467
$ws.Range("Q5").Value = 30
$ws.Range("R5").Value = "Hortaliza"

# Ensure the date cells use the same date/time number format as the other date column entries.
$ws.Range("D4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 8 appended at the end: week of 24-Aug-2021 (Elquí origin)
$ws.Range("A8").Value = 12
$ws.Range("B8").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C8").Value = "Metropolitana"
$ws.Range("D8").Value = 44432
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = 100112013
$ws.Range("G8").Value = "Alcachofa"
$ws.Range("H8").Value = "Española"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 25
$ws.Range("K8").Value = 14000
$ws.Range("L8").Value = 14000
$ws.Range("M8").Value = 14000
$ws.Range("N8").Value = "`$/caja 30 unidades"
$ws.Range("O8").Value = "Provincia del Elquí"
$ws.Range("P8").Value = 467
$ws.Range("Q8").Value = 30
$ws.Range("R8").Value = "Hortaliza"

$ws.Range("D8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
